# "Generate Report for Handback" - update the localization-status report
# after a handback run that failed for the f1accf87... source file in both
# the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# 1) The f1accf87... row's status flips from "Ready for handoff" to
#    "Handback transform failed" everywhere it is shown: the Overview
#    summary columns (zh-cn/de-de) and each language sheet's Status column.
$newStatus = "Handback transform failed"
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# 2) Widen the "Error Detail" column (P) on both language sheets so the
#    new failure message is readable.
$zhcn.Range("P1").ColumnWidth = 39.17
$dede.Range("P1").ColumnWidth = 39.17

# 3) Record the handback/handoff filename mismatch error for the failed
#    row (row 3) on each language sheet's "Error Detail" column.
$zhcn.Range("P3").Value = "Handback file name: q31wpcmq.1kr is different with handoff file name: f1accf87-46eb-437f-8e5c-dc11a701df30.1293a5faf9d50cdb002504960179651e907241d8.zh-cn."
$dede.Range("P3").Value = "Handback file name: q31wpcmq.1kr is different with handoff file name: f1accf87-46eb-437f-8e5c-dc11a701df30.1293a5faf9d50cdb002504960179651e907241d8.de-de."
